$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(7, 8).Value = 5282.75
$ws.Cells.Item(7, 9).Value = 2000
$ws.Cells.Item(7, 11).Value = 2000
$ws.Cells.Item(7, 13).Value = -1888

$ws.Cells.Item(14, 8).Value = 5282.75
$ws.Cells.Item(14, 9).Value = 2000
$ws.Cells.Item(14, 11).Value = 2000
$ws.Cells.Item(14, 13).Value = -1809

$ws.Cells.Item(29, 8).Value = 2630.3333
$ws.Cells.Item(29, 10).Value = 5166
$ws.Cells.Item(29, 12).Value = 15498
$ws.Cells.Item(29, 14).Value = -16060

$ws.Cells.Item(87, 8).Value = 59532.4
$ws.Cells.Item(87, 10).Value = 59532.4
$ws.Cells.Item(87, 12).Value = 59532.4
$ws.Cells.Item(87, 14).Value = -62028.4

$ws.Cells.Item(90, 8).Value = 59532.4
$ws.Cells.Item(90, 10).Value = 59532.4
$ws.Cells.Item(90, 12).Value = 178597.2
$ws.Cells.Item(90, 14).Value = -191077.2

$ws.Cells.Item(92, 8).Value = 1976.9
$ws.Cells.Item(92, 9).Value = 853.1667
$ws.Cells.Item(92, 10).Value = 3662.5
$ws.Cells.Item(92, 11).Value = 853.1667
$ws.Cells.Item(92, 12).Value = 3662.5
$ws.Cells.Item(92, 13).Value = 394.8333
$ws.Cells.Item(92, 14).Value = -6158.5

$ws.Cells.Item(113, 8).Value = 1948.75
$ws.Cells.Item(113, 9).Value = 1931.6666
$ws.Cells.Item(113, 11).Value = 1931.6666
$ws.Cells.Item(113, 13).Value = 1322.3334

$ws.Cells.Item(121, 8).Value = 0
$ws.Cells.Item(121, 10).Value = 0
$ws.Cells.Item(121, 12).Value = 0
$ws.Cells.Item(121, 14).ClearContents()

$ws.Cells.Item(137, 8).Value = 2512.7942
$ws.Cells.Item(137, 9).Value = 1942.8667
$ws.Cells.Item(137, 11).Value = 5828.6001
$ws.Cells.Item(137, 13).Value = -3278.6001

$ws.Cells.Item(141, 8).Value = 2365.6667
$ws.Cells.Item(141, 9).Value = 1818.8
$ws.Cells.Item(141, 11).Value = 5456.4
$ws.Cells.Item(141, 13).Value = -276.3999999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(23, 8).Value = 51669
$ws.Cells.Item(23, 10).Value = 51669
$ws.Cells.Item(23, 12).Value = 51669
$ws.Cells.Item(23, 14).Value = -52187

$ws.Cells.Item(46, 8).Value = 3999.3333
$ws.Cells.Item(46, 9).Value = 3499
$ws.Cells.Item(46, 11).Value = 3499
$ws.Cells.Item(46, 13).Value = -3180

$ws.Cells.Item(61, 8).Value = 4362.1665
$ws.Cells.Item(61, 9).Value = 3934.8
$ws.Cells.Item(61, 11).Value = 3934.8
$ws.Cells.Item(61, 13).Value = -3722.8

$ws.Cells.Item(62, 8).Value = 36666.332
$ws.Cells.Item(62, 10).Value = 36666.332
$ws.Cells.Item(62, 12).Value = 36666.332
$ws.Cells.Item(62, 14).Value = -37914.332

$ws.Cells.Item(65, 8).Value = 36666.332
$ws.Cells.Item(65, 10).Value = 36666.332
$ws.Cells.Item(65, 12).Value = 109998.996
$ws.Cells.Item(65, 14).Value = -116238.996

$ws.Cells.Item(68, 8).Value = 20000
$ws.Cells.Item(68, 10).Value = 20000
$ws.Cells.Item(68, 12).Value = 20000
$ws.Cells.Item(68, 14).Value = -21622

$ws.Cells.Item(71, 8).Value = 20000
$ws.Cells.Item(71, 10).Value = 20000
$ws.Cells.Item(71, 12).Value = 60000
$ws.Cells.Item(71, 14).Value = -68112

$ws.Cells.Item(74, 8).Value = 6012.6313
$ws.Cells.Item(74, 9).Value = 3232
$ws.Cells.Item(74, 10).Value = 9836
$ws.Cells.Item(74, 11).Value = 3232
$ws.Cells.Item(74, 12).Value = 9836
$ws.Cells.Item(74, 13).Value = -2358
$ws.Cells.Item(74, 14).Value = -11584

$ws.Cells.Item(76, 8).Value = 22000
$ws.Cells.Item(76, 10).Value = 22000
$ws.Cells.Item(76, 12).Value = 22000
$ws.Cells.Item(76, 14).Value = -22676

$ws.Cells.Item(77, 8).Value = 6012.6313
$ws.Cells.Item(77, 9).Value = 3232
$ws.Cells.Item(77, 10).Value = 9836
$ws.Cells.Item(77, 11).Value = 16160
$ws.Cells.Item(77, 12).Value = 49180
$ws.Cells.Item(77, 13).Value = -11792
$ws.Cells.Item(77, 14).Value = -57916

$ws.Cells.Item(79, 8).Value = 22000
$ws.Cells.Item(79, 10).Value = 22000
$ws.Cells.Item(79, 12).Value = 22000
$ws.Cells.Item(79, 14).Value = -24340

$ws.Cells.Item(132, 8).Value = 3860.5
$ws.Cells.Item(132, 9).Value = 1790.875
$ws.Cells.Item(132, 10).Value = 7999.75
$ws.Cells.Item(132, 11).Value = 5372.625
$ws.Cells.Item(132, 12).Value = 23999.25
$ws.Cells.Item(132, 13).Value = -2842.625
$ws.Cells.Item(132, 14).Value = -29059.25

$ws.Cells.Item(136, 8).Value = 4362.1665
$ws.Cells.Item(136, 9).Value = 3934.8
$ws.Cells.Item(136, 11).Value = 11804.4
$ws.Cells.Item(136, 13).Value = -9254.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1999.5
$ws.Cells.Item(94, 9).Value = 1999.5
$ws.Cells.Item(94, 11).Value = 1999.5
$ws.Cells.Item(94, 13).Value = -1548.5

$ws.Cells.Item(102, 8).Value = 7682.375
$ws.Cells.Item(102, 9).Value = 7682.375
$ws.Cells.Item(102, 11).Value = 7682.375
$ws.Cells.Item(102, 13).Value = -4437.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(102, 8).Value = 40120.5
$ws.Cells.Item(102, 10).Value = 40120.5
$ws.Cells.Item(102, 12).Value = 40120.5
$ws.Cells.Item(102, 14).Value = -44988.5

$ws.Cells.Item(105, 8).Value = 1246.1428
$ws.Cells.Item(105, 9).Value = 1246.1428
$ws.Cells.Item(105, 11).Value = 1246.1428
$ws.Cells.Item(105, 13).Value = 500.8571999999999

$ws.Cells.Item(132, 8).Value = 5211.28
$ws.Cells.Item(132, 9).Value = 4419.6
$ws.Cells.Item(132, 11).Value = 13258.8
$ws.Cells.Item(132, 13).Value = -10728.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 5624.875
$ws.Cells.Item(113, 10).Value = 5624.875
$ws.Cells.Item(113, 12).Value = 16874.625
$ws.Cells.Item(113, 14).Value = -21214.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 76.25
$ws.Cells.Item(2, 9).Value = 68.333336
$ws.Cells.Item(2, 11).Value = 68.333336
$ws.Cells.Item(2, 13).Value = 44.666664

$ws.Cells.Item(14, 8).Value = 6493.1113
$ws.Cells.Item(14, 9).Value = 101.6
$ws.Cells.Item(14, 10).Value = 14482.5
$ws.Cells.Item(14, 11).Value = 101.6
$ws.Cells.Item(14, 12).Value = 14482.5
$ws.Cells.Item(14, 13).Value = 66.40000000000001
$ws.Cells.Item(14, 14).Value = -14818.5

$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 14).ClearContents()

$ws.Cells.Item(46, 8).Value = 7963.6
$ws.Cells.Item(46, 9).Value = 5487
$ws.Cells.Item(46, 11).Value = 5487
$ws.Cells.Item(46, 13).Value = -5331

$ws.Cells.Item(55, 8).Value = 5051.9
$ws.Cells.Item(55, 9).Value = 6166
$ws.Cells.Item(55, 10).Value = 3937.8
$ws.Cells.Item(55, 11).Value = 6166
$ws.Cells.Item(55, 12).Value = 3937.8
$ws.Cells.Item(55, 13).Value = -5839
$ws.Cells.Item(55, 14).Value = -4591.8

$ws.Cells.Item(70, 8).Value = 6140.857
$ws.Cells.Item(70, 9).Value = 5497.8335
$ws.Cells.Item(70, 11).Value = 5497.8335
$ws.Cells.Item(70, 13).Value = -5227.8335

$ws.Cells.Item(73, 8).Value = 6140.857
$ws.Cells.Item(73, 9).Value = 5497.8335
$ws.Cells.Item(73, 11).Value = 5497.8335
$ws.Cells.Item(73, 13).Value = -4561.8335

$ws.Cells.Item(99, 8).Value = 7617.5
$ws.Cells.Item(99, 9).Value = 4823.3335
$ws.Cells.Item(99, 11).Value = 4823.3335
$ws.Cells.Item(99, 13).Value = -2577.3335

$ws.Cells.Item(113, 8).Value = 9874.25
$ws.Cells.Item(113, 9).Value = 8999
$ws.Cells.Item(113, 11).Value = 8999
$ws.Cells.Item(113, 13).Value = -6829

$ws.Cells.Item(132, 8).Value = 59588.39
$ws.Cells.Item(132, 9).Value = 62505.35
$ws.Cells.Item(132, 10).Value = 10000
$ws.Cells.Item(132, 11).Value = 187516.05
$ws.Cells.Item(132, 12).Value = 30000
$ws.Cells.Item(132, 13).Value = -184986.05
$ws.Cells.Item(132, 14).Value = -35060

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 4254.6665
$ws.Cells.Item(22, 9).Value = 1000
$ws.Cells.Item(22, 10).Value = 4905.6
$ws.Cells.Item(22, 11).Value = 1000
$ws.Cells.Item(22, 12).Value = 4905.6
$ws.Cells.Item(22, 13).Value = -705
$ws.Cells.Item(22, 14).Value = -5495.6

$ws.Cells.Item(27, 8).Value = 4254.6665
$ws.Cells.Item(27, 9).Value = 1000
$ws.Cells.Item(27, 10).Value = 4905.6
$ws.Cells.Item(27, 11).Value = 1000
$ws.Cells.Item(27, 12).Value = 4905.6
$ws.Cells.Item(27, 13).Value = -893
$ws.Cells.Item(27, 14).Value = -5119.6

$ws.Cells.Item(47, 8).Value = 39666.668
$ws.Cells.Item(47, 10).Value = 44500
$ws.Cells.Item(47, 12).Value = 44500
$ws.Cells.Item(47, 14).Value = -45480

$ws.Cells.Item(52, 8).Value = 39666.668
$ws.Cells.Item(52, 10).Value = 44500
$ws.Cells.Item(52, 12).Value = 44500
$ws.Cells.Item(52, 14).Value = -44966

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(12, 8).Value = 500
$ws.Cells.Item(12, 9).Value = 500
$ws.Cells.Item(12, 11).Value = 500
$ws.Cells.Item(12, 13).Value = -358

$ws.Cells.Item(17, 8).Value = 2501.6
$ws.Cells.Item(17, 9).Value = 3077
$ws.Cells.Item(17, 10).Value = 200
$ws.Cells.Item(17, 11).Value = 3077
$ws.Cells.Item(17, 12).Value = 200
$ws.Cells.Item(17, 13).Value = -2905
$ws.Cells.Item(17, 14).Value = -544

$ws.Cells.Item(113, 8).Value = 745.3333
$ws.Cells.Item(113, 9).Value = 1056.5
$ws.Cells.Item(113, 11).Value = 3169.5
$ws.Cells.Item(113, 13).Value = -999.5

$ws.Cells.Item(126, 8).Value = 3559.7144
$ws.Cells.Item(126, 9).Value = 1259
$ws.Cells.Item(126, 11).Value = 3777
$ws.Cells.Item(126, 13).Value = -1307
